$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.102082522762076
$ws.Range("D2").Value = 1.094257490690614
$ws.Range("E2").Value = 1.113391606586077
$ws.Range("F2").Value = 1.114997830539742
$ws.Range("I2").Value = 1.059956023969702
$ws.Range("J2").Value = 1.106857096876727
$ws.Range("K2").Value = 1.096895294349004
$ws.Range("L2").Value = 1.115981648082747
$ws.Range("M2").Value = 1.11758394284288
$ws.Range("N2").Value = 1.108428960912575
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.10422154200374
$ws.Range("D3").Value = 1.095949256407573
$ws.Range("E3").Value = 1.115464572156477
$ws.Range("F3").Value = 1.11700462431261
$ws.Range("I3").Value = 1.060581209124672
$ws.Range("J3").Value = 1.108660130122111
$ws.Range("K3").Value = 1.098405816743618
$ws.Range("L3").Value = 1.117876004065712
$ws.Range("M3").Value = 1.119412566801111
$ws.Range("N3").Value = 1.110234554671979
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.105599835991208
$ws.Range("D4").Value = 1.097038368185687
$ws.Range("E4").Value = 1.1168004306016
$ws.Range("F4").Value = 1.118297647980167
$ws.Range("I4").Value = 1.060981636179896
$ws.Range("J4").Value = 1.109820776565464
$ws.Range("K4").Value = 1.0993771794384
$ws.Range("L4").Value = 1.119095881899714
$ws.Range("M4").Value = 1.120589880915466
$ws.Range("N4").Value = 1.11139684936641
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.10617791521561
$ws.Range("D5").Value = 1.097494921465341
$ws.Range("E5").Value = 1.117360741291013
$ws.Range("F5").Value = 1.118839945891896
$ws.Range("I5").Value = 1.061149003011765
$ws.Range("J5").Value = 1.110307294651117
$ws.Range("K5").Value = 1.099784116697289
$ws.Range("L5").Value = 1.119607335066373
$ws.Range("M5").Value = 1.121083432109933
$ws.Range("N5").Value = 1.111884058363549
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.106274898684202
$ws.Range("D6").Value = 1.097571502715146
$ws.Range("E6").Value = 1.117454745477548
$ws.Range("F6").Value = 1.118930925331829
$ws.Range("I6").Value = 1.06117704792098
$ws.Range("J6").Value = 1.110388900839612
$ws.Range("K6").Value = 1.099852360461777
$ws.Range("L6").Value = 1.119693129993999
$ws.Range("M6").Value = 1.12116622071361
$ws.Range("N6").Value = 1.111965780442192
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.105607565598013
$ws.Range("D7").Value = 1.097044473783927
$ws.Range("E7").Value = 1.116807922506335
$ws.Range("F7").Value = 1.118304899218307
$ws.Range("I7").Value = 1.060983876350748
$ws.Range("J7").Value = 1.109827282977454
$ws.Range("K7").Value = 1.099382622515685
$ws.Range("L7").Value = 1.11910272135998
$ws.Range("M7").Value = 1.12059648120198
$ws.Range("N7").Value = 1.11140336501825
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.102806633291052
$ws.Range("D8").Value = 1.094830400105281
$ws.Range("E8").Value = 1.114093329604487
$ws.Range("F8").Value = 1.115677193493949
$ws.Range("I8").Value = 1.060168165876118
$ws.Range("J8").Value = 1.107467706993061
$ws.Range("K8").Value = 1.097407048860958
$ws.Range("L8").Value = 1.116623092760034
$ws.Range("M8").Value = 1.118203178927321
$ws.Range("N8").Value = 1.109040438165313
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.097825220251417
$ws.Range("D9").Value = 1.090885095374592
$ws.Range("E9").Value = 1.109266479333308
$ws.Range("F9").Value = 1.111003359379538
$ws.Range("I9").Value = 1.058698800623202
$ws.Range("J9").Value = 1.103262353031657
$ws.Range("K9").Value = 1.093878445989734
$ws.Range("L9").Value = 1.112207222124352
$ws.Range("M9").Value = 1.113939228856248
$ws.Range("N9").Value = 1.104829112118922
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.094471455128121
$ws.Range("D10").Value = 1.088223825177023
$ws.Range("E10").Value = 1.106017498461618
$ws.Range("F10").Value = 1.107856387145977
$ws.Range("I10").Value = 1.057697009426221
$ws.Range("J10").Value = 1.100425103254182
$ws.Range("K10").Value = 1.091492654907747
$ws.Range("L10").Value = 1.109230269500414
$ws.Range("M10").Value = 1.111063476907644
$ws.Range("N10").Value = 1.101987833121326
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.093010997907081
$ws.Range("D11").Value = 1.087063736425201
$ws.Range("E11").Value = 1.104602852858471
$ws.Range("F11").Value = 1.10648592501369
$ws.Range("I11").Value = 1.057257792958209
$ws.Range("J11").Value = 1.099188158614302
$ws.Range("K11").Value = 1.090451317775457
$ws.Range("L11").Value = 1.107932972532389
$ws.Range("M11").Value = 1.109809993944096
$ws.Range("N11").Value = 1.100749131878178
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.092467239402929
$ws.Range("D12").Value = 1.086631633434177
$ws.Range("E12").Value = 1.104076179050462
$ws.Range("F12").Value = 1.105975666008234
$ws.Range("I12").Value = 1.057093817946748
$ws.Range("J12").Value = 1.098727406409209
$ws.Range("K12").Value = 1.090063245856778
$ws.Range("L12").Value = 1.107449822802638
$ws.Range("M12").Value = 1.109343118729661
$ws.Range("N12").Value = 1.100287725352105
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.092583935996385
$ws.Range("D13").Value = 1.086724375549595
$ws.Range("E13").Value = 1.104189207795486
$ws.Range("F13").Value = 1.10608517356399
$ws.Range("I13").Value = 1.057129028951095
$ws.Range("J13").Value = 1.098826298521141
$ws.Range("K13").Value = 1.090146546689328
$ws.Range("L13").Value = 1.107553518370133
$ws.Range("M13").Value = 1.109443323339733
$ws.Range("N13").Value = 1.100386757902179
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.092966076995102
$ws.Range("D14").Value = 1.087028043205377
$ws.Range("E14").Value = 1.104559342768157
$ws.Range("F14").Value = 1.106443771706336
$ws.Range("I14").Value = 1.05724425575574
$ws.Range("J14").Value = 1.099150099330059
$ws.Range("K14").Value = 1.090419265772016
$ws.Range("L14").Value = 1.1078930614593
$ws.Range("M14").Value = 1.109771428120628
$ws.Range("N14").Value = 1.100711018545387
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.09320135591592
$ws.Range("D15").Value = 1.087214983656861
$ws.Range("E15").Value = 1.104787233486998
$ws.Range("F15").Value = 1.106664554647465
$ws.Range("I15").Value = 1.057315140325377
$ws.Range("J15").Value = 1.099349430938591
$ws.Range("K15").Value = 1.090587127419414
$ws.Range("L15").Value = 1.108102095028375
$ws.Range("M15").Value = 1.109973414207248
$ws.Range("N15").Value = 1.100910633227667
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.094568203590874
$ws.Range("D16").Value = 1.088300650685287
$ws.Range("E16").Value = 1.106111216002869
$ws.Range("F16").Value = 1.107947172747124
$ws.Range("I16").Value = 1.057726043091716
$ws.Range("J16").Value = 1.100507015365614
$ws.Range("K16").Value = 1.091561588163654
$ws.Range("L16").Value = 1.109316189859766
$ws.Range("M16").Value = 1.111146489425066
$ws.Range("N16").Value = 1.102069861557352
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.095423353326912
$ws.Range("D17").Value = 1.088979566430472
$ws.Range("E17").Value = 1.106939596734186
$ws.Range("F17").Value = 1.108749610962443
$ws.Range("I17").Value = 1.057982326812336
$ws.Range("J17").Value = 1.101230865361912
$ws.Range("K17").Value = 1.092170605560142
$ws.Range("L17").Value = 1.110075524097538
$ws.Range("M17").Value = 1.111880092108312
$ws.Range("N17").Value = 1.102794739503672
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.095921353931989
$ws.Range("D18").Value = 1.08937482172956
$ws.Range("E18").Value = 1.107422025670823
$ws.Range("F18").Value = 1.109216909250881
$ws.Range("I18").Value = 1.058131289402034
$ws.Range("J18").Value = 1.101652267003186
$ws.Range("K18").Value = 1.092525038841742
$ws.Range("L18").Value = 1.110517636666069
$ws.Range("M18").Value = 1.112307195193635
$ws.Range("N18").Value = 1.103216739583612
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.096091025816774
$ws.Range("D19").Value = 1.089509468187243
$ws.Range("E19").Value = 1.107586395117061
$ws.Range("F19").Value = 1.109376119703653
$ws.Range("I19").Value = 1.058181993443804
$ws.Range("J19").Value = 1.101795818084917
$ws.Range("K19").Value = 1.092645757268454
$ws.Range("L19").Value = 1.110668251970586
$ws.Range("M19").Value = 1.112452692502178
$ws.Range("N19").Value = 1.103360494524344
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.095331686240197
$ws.Range("D20").Value = 1.088906802456478
$ws.Range("E20").Value = 1.106850797323439
$ws.Range("F20").Value = 1.108663594807698
$ws.Range("I20").Value = 1.057954884229803
$ws.Range("J20").Value = 1.101153286910922
$ws.Range("K20").Value = 1.092105346333146
$ws.Range("L20").Value = 1.109994137126924
$ws.Range("M20").Value = 1.111801465955281
$ws.Range("N20").Value = 1.102717050882383
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.092853581641004
$ws.Range("D21").Value = 1.086938653874523
$ws.Range("E21").Value = 1.104450380988269
$ws.Range("F21").Value = 1.106338207140828
$ws.Range("I21").Value = 1.0572103473723
$ws.Range("J21").Value = 1.099054784114259
$ws.Range("K21").Value = 1.090338992181549
$ws.Range("L21").Value = 1.107793109989914
$ws.Range("M21").Value = 1.1096748449384
$ws.Range("N21").Value = 1.100615567971047
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.091288075125885
$ws.Range("D22").Value = 1.085694274532253
$ws.Range("E22").Value = 1.102934116001179
$ws.Range("F22").Value = 1.104869134132173
$ws.Range("I22").Value = 1.056737415591575
$ws.Range("J22").Value = 1.097727857033112
$ws.Range("K22").Value = 1.089221035185358
$ws.Range("L22").Value = 1.106401837090256
$ws.Range("M22").Value = 1.108330354891539
$ws.Range("N22").Value = 1.099286756501245
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.092118696797409
$ws.Range("D23").Value = 1.086354610699271
$ws.Range("E23").Value = 1.103738595537313
$ws.Range("F23").Value = 1.105648594173227
$ws.Range("I23").Value = 1.05698858653656
$ws.Range("J23").Value = 1.098432010159148
$ws.Range("K23").Value = 1.089814394998319
$ws.Range("L23").Value = 1.107140090562613
$ws.Range("M23").Value = 1.109043807467552
$ws.Range("N23").Value = 1.099991909605487
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.095373109102092
$ws.Range("D24").Value = 1.088939683664357
$ws.Range("E24").Value = 1.106890924275604
$ws.Range("F24").Value = 1.108702464122674
$ws.Range("I24").Value = 1.057967285968246
$ws.Range("J24").Value = 1.101188343768382
$ws.Range("K24").Value = 1.092134836628033
$ws.Range("L24").Value = 1.110030914848507
$ws.Range("M24").Value = 1.11183699618445
$ws.Range("N24").Value = 1.102752157524602
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.099118670036773
$ws.Range("D25").Value = 1.091910407926748
$ws.Range("E25").Value = 1.110519672306737
$ws.Range("F25").Value = 1.112216998725465
$ws.Range("I25").Value = 1.059082529150681
$ws.Range("J25").Value = 1.104355343621941
$ws.Range("K25").Value = 1.094796445997215
$ws.Range("L25").Value = 1.113354519386903
$ws.Range("M25").Value = 1.11504726899753
$ws.Range("N25").Value = 1.105923654881213
